$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E10").Value = 424

$ws.Range("E12").Value = 421

$ws.Range("E13").Value = 111

$ws.Range("E15").Value = 139

$ws.Range("E16").Value = 175

$ws.Range("E25").Value = 225

$ws.Range("E29").Value = 149

$ws.Range("E30").Value = 186
$ws.Range("F30").Value = 107
$ws.Range("H30").Value = 107

$ws.Range("E32").Value = 165

$ws.Range("E33").Value = 253

$ws.Range("E34").Value = 187

$ws.Range("E36").Value = 60
$ws.Range("F36").Value = 37
$ws.Range("H36").Value = 37

$ws.Range("E40").Value = 232

$ws.Range("E41").Value = 343
$ws.Range("F41").Value = 162
$ws.Range("H41").Value = 162

$ws.Range("E42").Value = 315

$ws.Range("E44").Value = 269

$ws.Range("E47").Value = 386
$ws.Range("F47").Value = 189
$ws.Range("H47").Value = 189

$wb.Save()
